$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column whose new value still looks like a plain
# number (e.g. "1.007", "8.210") must be pre-formatted as Text so Excel
# keeps the exact source string (incl. trailing zeros) instead of silently
# parsing it into a number and dropping the formatting.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.385.78'
$ws.Range("E2").Value = '  +1.33%  '

$ws.Range("D3").Value = '1.686.34'
$ws.Range("E3").Value = '  +1.11%  '

$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").Value = '218.57'
$ws.Range("E5").Value = '  +0.88%  '

$ws.Range("D6").Value = '0.5599'
$ws.Range("E6").Value = '  +9.92%  '

$ws.Range("E7").Value = '  +0.24%  '

$ws.Range("D8").Value = '0.2716'
$ws.Range("E8").Value = '  +2.25%  '

$ws.Range("D9").Value = '0.06524'
$ws.Range("E9").Value = '  +2.12%  '

$ws.Range("D10").Value = '22.18'
$ws.Range("E10").Value = '  +1.96%  '

$ws.Range("D11").Value = '0.07573'
$ws.Range("E11").Value = '  +1.72%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.689.24'
$ws.Range("E12").Value = '  +1.09%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '4.557'
$ws.Range("E13").Value = '  +1.04%  '

$ws.Range("D14").Value = '0.5826'
$ws.Range("E14").Value = '  +0.04%  '

$ws.Range("D15").Value = '0.000008484'
$ws.Range("E15").Value = '  -0.68%  '

$ws.Range("D16").Value = '65.57'
$ws.Range("E16").Value = '  +2.03%  '

$ws.Range("D17").Value = '26.399.48'
$ws.Range("E17").Value = '  +1.11%  '

$ws.Range("D18").Value = '4.954'
$ws.Range("E18").Value = '  +0.44%  '

$ws.Range("E19").Value = '  +0.25%  '

$ws.Range("D20").Value = '10.95'
$ws.Range("E20").Value = '  +1.62%  '

$ws.Range("D21").Value = '191.74'
$ws.Range("E21").Value = '  +0.10%  '

$ws.Range("D22").Value = '6.253'
$ws.Range("E22").Value = '  +0.95%  '

$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("D24").Value = '148.72'
$ws.Range("E24").Value = '  +2.74%  '

$ws.Range("E25").Value = '  +11.45%  '

$ws.Range("D26").Value = '7.922'
$ws.Range("E26").Value = '  +4.16%  '

$ws.Range("D27").Value = '15.88'
$ws.Range("E27").Value = '  +1.44%  '

$ws.Range("D28").Value = '0.06313'
$ws.Range("E28").Value = '  -3.94%  '

$ws.Range("E29").Value = '  +4.06%  '

$ws.Range("D30").Value = '1.323'
$ws.Range("E30").Value = '  +0.67%  '

$ws.Range("D31").Value = '3.605'
$ws.Range("E31").Value = '  +1.86%  '

$ws.Range("D32").Value = '3.588'
$ws.Range("E32").Value = '  +2.27%  '

$ws.Range("D33").Value = '1.676'
$ws.Range("E33").Value = '  +1.39%  '

$ws.Range("E34").Value = '  +2.55%  '

$ws.Range("E35").Value = '  +1.94%  '

$ws.Range("D36").Value = '2.401'

$ws.Range("D37").Value = '2.716'
$ws.Range("E37").Value = '  +1.19%  '

$ws.Range("D38").Value = '6.243'
$ws.Range("E38").Value = '  -0.73%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01636'
$ws.Range("E39").Value = '  +2.37%  '

$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.116.68'
$ws.Range("E40").Value = '  +2.39%  '

$ws.Range("D41").Value = '0.8763'
$ws.Range("E41").Value = '  +0.55%  '

$ws.Range("E42").Value = '  +0.45%  '

$ws.Range("D43").Value = '100.69'
$ws.Range("E43").Value = '  -0.43%  '

$ws.Range("D44").Value = '1.835.74'
$ws.Range("E44").Value = '  +1.09%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '57.47'
$ws.Range("E45").Value = '  +1.96%  '

$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.00000000108'
$ws.Range("E46").Value = '  -2.76%  '

$ws.Range("D47").Value = '8.210'
$ws.Range("E47").Value = '  +1.63%  '

$ws.Range("D48").Value = '1.006'
$ws.Range("E48").Value = '  -0.25%  '

$ws.Range("D49").Value = '0.05280'
$ws.Range("E49").Value = '  +0.91%  '

$ws.Range("D50").Value = '6.097'
$ws.Range("E50").Value = '  +0.99%  '

$ws.Range("D51").Value = '0.4298'
$ws.Range("E51").Value = '  +0.24%  '
